$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 297.0263
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 290.32
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 870.96
$ws.Range("M17").Value = -2232
$ws.Range("N17").Value = -1206.96
$ws.Range("H92").Value = 806.4375
$ws.Range("I92").Value = 323.08694
$ws.Range("J92").Value = 2041.6666
$ws.Range("K92").Value = 323.08694
$ws.Range("L92").Value = 2041.6666
$ws.Range("M92").Value = 924.91306
$ws.Range("N92").Value = -4537.6666
$ws.Range("H100").Value = 57210.5
$ws.Range("I100").Value = 101450.9
$ws.Range("J100").Value = 1910
$ws.Range("K100").Value = 101450.9
$ws.Range("L100").Value = 1910
$ws.Range("M100").Value = -100909.9
$ws.Range("N100").Value = -2992
$ws.Range("H111").Value = 1613.1111
$ws.Range("I111").Value = 1598.75
$ws.Range("J111").Value = 1624.6
$ws.Range("K111").Value = 4796.25
$ws.Range("L111").Value = 4873.799999999999
$ws.Range("M111").Value = -1729.25
$ws.Range("N111").Value = -11007.8
$ws.Range("H137").Value = 311864.53
$ws.Range("I137").Value = 351837.1
$ws.Range("J137").Value = 85353.336
$ws.Range("K137").Value = 1055511.3
$ws.Range("L137").Value = 256060.008
$ws.Range("M137").Value = -1052961.3
$ws.Range("N137").Value = -261160.008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 816288.8
$ws.Range("I32").Value = 5891.423
$ws.Range("K32").Value = 5891.423
$ws.Range("M32").Value = -5604.423
$ws.Range("H97").Value = 1318.3684
$ws.Range("I97").Value = 659
$ws.Range("J97").Value = 2225
$ws.Range("K97").Value = 659
$ws.Range("L97").Value = 2225
$ws.Range("M97").Value = -163
$ws.Range("N97").Value = -3217
$ws.Range("H122").Value = 6453214.5
$ws.Range("I122").Value = 10001200
$ws.Range("K122").Value = 30003600
$ws.Range("M122").Value = -30001150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5988.36
$ws.Range("I20").Value = 7783.5884
$ws.Range("J20").Value = 2173.5
$ws.Range("K20").Value = 7783.5884
$ws.Range("L20").Value = 2173.5
$ws.Range("M20").Value = -7536.5884
$ws.Range("N20").Value = -2667.5
$ws.Range("H86").Value = 1529.4375
$ws.Range("I86").Value = 1555.3334
$ws.Range("J86").Value = 1451.75
$ws.Range("K86").Value = 1555.3334
$ws.Range("L86").Value = 1451.75
$ws.Range("M86").Value = -432.3334
$ws.Range("N86").Value = -3697.75
$ws.Range("H89").Value = 1529.4375
$ws.Range("I89").Value = 1555.3334
$ws.Range("J89").Value = 1451.75
$ws.Range("K89").Value = 7776.666999999999
$ws.Range("L89").Value = 7258.75
$ws.Range("M89").Value = -2160.666999999999
$ws.Range("N89").Value = -18490.75
$ws.Range("H94").Value = 768.36365
$ws.Range("I94").Value = 584.88464
$ws.Range("J94").Value = 1449.8572
$ws.Range("K94").Value = 584.88464
$ws.Range("L94").Value = 1449.8572
$ws.Range("M94").Value = -133.88464
$ws.Range("N94").Value = -2351.8572
$ws.Range("H99").Value = 1107.5385
$ws.Range("I99").Value = 969.9
$ws.Range("K99").Value = 969.9
$ws.Range("M99").Value = 528.1
$ws.Range("H105").Value = 2045.7142
$ws.Range("I105").Value = 1848.8889
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 1848.8889
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -101.8888999999999
$ws.Range("N105").Value = -5894

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26176.1
$ws.Range("I31").Value = 9292.333000000001
$ws.Range("J31").Value = 51501.75
$ws.Range("K31").Value = 9292.333000000001
$ws.Range("L31").Value = 51501.75
$ws.Range("M31").Value = -8997.333000000001
$ws.Range("N31").Value = -52091.75
$ws.Range("H34").Value = 26176.1
$ws.Range("I34").Value = 9292.333000000001
$ws.Range("J34").Value = 51501.75
$ws.Range("K34").Value = 9292.333000000001
$ws.Range("L34").Value = 51501.75
$ws.Range("M34").Value = -9090.333000000001
$ws.Range("N34").Value = -51905.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 100549.8
$ws.Range("J5").Value = 250361
$ws.Range("L5").Value = 751083
$ws.Range("N5").Value = -751307
$ws.Range("H87").Value = 2338
$ws.Range("I87").Value = 2338
$ws.Range("K87").Value = 7014
$ws.Range("M87").Value = -5766
$ws.Range("H90").Value = 2338
$ws.Range("I90").Value = 2338
$ws.Range("K90").Value = 21042
$ws.Range("M90").Value = -14802
$ws.Range("H117").Value = 876.3333
$ws.Range("H135").Value = 100549.8
$ws.Range("J135").Value = 250361
$ws.Range("L135").Value = 2253249
$ws.Range("N135").Value = -2258319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26803.244
$ws.Range("I70").Value = 33836.65
$ws.Range("K70").Value = 33836.65
$ws.Range("M70").Value = -33566.65
$ws.Range("H73").Value = 26803.244
$ws.Range("I73").Value = 33836.65
$ws.Range("K73").Value = 33836.65
$ws.Range("M73").Value = -32900.65
$ws.Range("H97").Value = 50001096
$ws.Range("I97").Value = 911.2857
$ws.Range("J97").Value = 76924270
$ws.Range("K97").Value = 911.2857
$ws.Range("L97").Value = 76924270
$ws.Range("M97").Value = -415.2857
$ws.Range("N97").Value = -76925262
$ws.Range("H122").Value = 1238
$ws.Range("I122").Value = 1013
$ws.Range("J122").Value = 1500.5
$ws.Range("K122").Value = 3039
$ws.Range("L122").Value = 4501.5
$ws.Range("M122").Value = -589
$ws.Range("N122").Value = -9401.5
$ws.Range("H132").Value = 30305030
$ws.Range("I132").Value = 43479410
$ws.Range("J132").Value = 3959.5
$ws.Range("K132").Value = 130438230
$ws.Range("L132").Value = 11878.5
$ws.Range("M132").Value = -130435700
$ws.Range("N132").Value = -16938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3386.2354
$ws.Range("I82").Value = 1295.5555
$ws.Range("J82").Value = 5738.25
$ws.Range("K82").Value = 1295.5555
$ws.Range("L82").Value = 5738.25
$ws.Range("M82").Value = -934.5554999999999
$ws.Range("N82").Value = -6460.25
$ws.Range("H85").Value = 3386.2354
$ws.Range("I85").Value = 1295.5555
$ws.Range("J85").Value = 5738.25
$ws.Range("K85").Value = 1295.5555
$ws.Range("L85").Value = 5738.25
$ws.Range("M85").Value = -47.55549999999994
$ws.Range("N85").Value = -8234.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5195.4346
$ws.Range("I122").Value = 9876.182000000001
$ws.Range("J122").Value = 904.75
$ws.Range("K122").Value = 29628.546
$ws.Range("L122").Value = 2714.25
$ws.Range("M122").Value = -27178.546
$ws.Range("N122").Value = -7614.25
